$d = $word.ActiveDocument

# --- Change 1: paragraph 9 'Pagination der Kunden-, ...' -> split with proofErr ---
$p9 = $d.Paragraphs.Item(9)
$p9.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Pagination</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> der Kunden-, Events- und Transaktionslisten inklusive der Kunden- und Transaktionssuche</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# --- Change 2: paragraph 25 'Nicht lauffaehiger Code ...' -> split with proofErr ---
$p25 = $d.Paragraphs.Item(25)
$p25.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t>Nicht lauffähiger Code nach dem Sprintende am master-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>branch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ist nicht akzeptabel</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# --- Change 3: paragraphs 26+27 'Massnahme: ...' + '3. Sprint' heading -> split runs, move bookmark ---
$p26 = $d.Paragraphs.Item(26)
$p27 = $d.Paragraphs.Item(27)
$range26_27 = $d.Range($p26.Range.Start, $p27.Range.End)
$range26_27.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t>Maßnahme: bevor der Code dann wirklich auf den master-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>branch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gepushed</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> wird, sollte eine Person den Code nochmal außerhalb von </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Intellij</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> und jeder Entwicklungsumgebung nur über die Command-Line manuell kompilieren und testen.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>. Sprint</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# --- Change 4: paragraph 28 (empty) + new paragraph -> 'realistische Testdaten' + new retro item ---
$p28 = $d.Paragraphs.Item(28)
$p28.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t>realistische Testdaten</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">wir haben das Problem, dass wir vom Programm aus nicht die automatisch generierten </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LastModifiedAt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> und </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CreatedAt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Timestamps</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in der Datenbank zu manipulieren, was ja grundsätzlich gut ist, aber für die Tests nicht sehr hilfreich. Als Maßnahme könnten wir ein SQL Script schreiben, das uns alle gewünschten Daten direkt in die Datenbank einfügt falls so etwas überhaupt in unserer </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Hibernate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> / Spring Einstellung möglich ist.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

